$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear the old (now obsolete) encryption-table values from row 2 ---
# (O2,P2,R2,U2,V2,W2,X2 held the q/w/e/t/y/u/i letters; the cipher table
#  row that used to live at row 2 is being replaced by a new key-row 6)
$ws.Range("O2").Value = $null
$ws.Range("P2").Value = $null
$ws.Range("R2").Value = $null
$ws.Range("U2").Value = $null
$ws.Range("V2").Value = $null
$ws.Range("W2").Value = $null
$ws.Range("X2").Value = $null

# --- Build the new row 6 "key 1" entries (N6:X6) ---
# N6 gets the "Uwaga" formatting used by the other key-number cells (N3:N5)
$ws.Range("N3").Copy()
$ws.Range("N6").PasteSpecial(-4122)
$ws.Range("N6").Value = 1

# O6:V6 get the "Dane wejściowe" formatting used by the rest of the table
$ws.Range("O2:V2").Copy()
$ws.Range("O6:V6").PasteSpecial(-4122)

# W6:X6 also pick up the same formatting but stay blank
$ws.Range("O2").Copy()
$ws.Range("W6:X6").PasteSpecial(-4122)

$ws.Range("O6").Value = "w"
$ws.Range("P6").Value = "[space]"
$ws.Range("Q6").Value = "q"
$ws.Range("R6").Value = "e"
$ws.Range("S6").Value = "t"
$ws.Range("T6").Value = "y"
$ws.Range("U6").Value = "u"
$ws.Range("V6").Value = "i"

$excel.CutCopyMode = 0

# --- Update the active selection to match the saved view state ---
$ws.Range("S9").Select()

$wb.Save()
